$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.665.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "'1.597.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'211.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "'0.515"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "'1.822.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "'1.594.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "'65.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "'26.654.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D19").Value = "'209.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("E21").Value = "  +3.91%  "
$ws.Range("D22").Value = "'4.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("D23").Value = "'2.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("D24").Value = "'8.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").Value = "'144.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "'7.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").Value = "'15.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "'0.0517"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("D31").Value = "'1.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("D34").Value = "'1.287.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("E35").Value = "  -6.83%  "
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("E40").Value = "  +18.73%  "
$ws.Range("D41").Value = "'5.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.33%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D44").Value = "'63.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").Value = "'1.734.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'90.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("D48").Value = "'0.101"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("D49").Value = "'0.0509"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "'7.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.26%  "
